$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 74
$ws.Range("H74").Value = 3400.15
$ws.Range("I74").Value = 3083.5833
$ws.Range("J74").Value = 3875
$ws.Range("K74").Value = 3083.5833
$ws.Range("L74").Value = 3875
$ws.Range("M74").Value = -2147.5833
$ws.Range("N74").Value = -5747

# Row 77
$ws.Range("H77").Value = 3400.15
$ws.Range("I77").Value = 3083.5833
$ws.Range("J77").Value = 3875
$ws.Range("K77").Value = 15417.9165
$ws.Range("L77").Value = 19375
$ws.Range("M77").Value = -10737.9165
$ws.Range("N77").Value = -28735

# Row 123
$ws.Range("H123").Value = 37700
$ws.Range("J123").Value = 37700
$ws.Range("L123").Value = 37700
$ws.Range("N123").Value = -47500

# Row 124
$ws.Range("H124").Value = 48022.2
$ws.Range("J124").Value = 48022.2
$ws.Range("L124").Value = 48022.2
$ws.Range("N124").Value = -57842.2

# Row 126
$ws.Range("H126").Value = 38614.168
$ws.Range("J126").Value = 38614.168
$ws.Range("L126").Value = 38614.168
$ws.Range("N126").Value = -48494.168

# Row 128
$ws.Range("H128").Value = 39829
$ws.Range("J128").Value = 39829
$ws.Range("L128").Value = 39829
$ws.Range("N128").Value = -49789


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 107
$ws.Range("H107").Value = 41111.5
$ws.Range("J107").Value = 41111.5
$ws.Range("L107").Value = 41111.5
$ws.Range("N107").Value = -48791.5

# Row 109
$ws.Range("H109").Value = 41934.6
$ws.Range("J109").Value = 41934.6
$ws.Range("L109").Value = 41934.6
$ws.Range("N109").Value = -44708.6

# Row 111
$ws.Range("H111").Value = 49616
$ws.Range("J111").Value = 49616
$ws.Range("L111").Value = 49616
$ws.Range("N111").Value = -57796

# Row 113
$ws.Range("H113").Value = 46619
$ws.Range("J113").Value = 46619
$ws.Range("L113").Value = 46619
$ws.Range("N113").Value = -55297

# Row 114
$ws.Range("H114").Value = 45970.668
$ws.Range("J114").Value = 45970.668
$ws.Range("L114").Value = 45970.668
$ws.Range("N114").Value = -54648.668


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 1773.1786
$ws.Range("I86").Value = 1825.2778
$ws.Range("J86").Value = 1679.4
$ws.Range("K86").Value = 1825.2778
$ws.Range("L86").Value = 1679.4
$ws.Range("M86").Value = -702.2778000000001
$ws.Range("N86").Value = -3925.4

# Row 88
$ws.Range("H88").Value = 26687.166
$ws.Range("J88").Value = 26687.166
$ws.Range("L88").Value = 26687.166
$ws.Range("N88").Value = -27499.166

# Row 89
$ws.Range("H89").Value = 1773.1786
$ws.Range("I89").Value = 1825.2778
$ws.Range("J89").Value = 1679.4
$ws.Range("K89").Value = 9126.389000000001
$ws.Range("L89").Value = 8397
$ws.Range("M89").Value = -3510.389000000001
$ws.Range("N89").Value = -19629

# Row 91
$ws.Range("H91").Value = 26687.166
$ws.Range("J91").Value = 26687.166
$ws.Range("L91").Value = 26687.166
$ws.Range("N91").Value = -29495.166

# Row 94
$ws.Range("H94").Value = 675.5263
$ws.Range("I94").Value = 664.0625
$ws.Range("J94").Value = 736.6667
$ws.Range("K94").Value = 664.0625
$ws.Range("L94").Value = 736.6667
$ws.Range("M94").Value = -213.0625
$ws.Range("N94").Value = -1638.6667

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 110
$ws.Range("H110").Value = 46900
$ws.Range("J110").Value = 46900
$ws.Range("L110").Value = 46900
$ws.Range("N110").Value = -55080

# Row 111
$ws.Range("H111").Value = 48702
$ws.Range("J111").Value = 48702
$ws.Range("L111").Value = 48702
$ws.Range("N111").Value = -56882

# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 133
$ws.Range("H133").Value = 27239.357
$ws.Range("J133").Value = 27239.357
$ws.Range("L133").Value = 27239.357
$ws.Range("N133").Value = -32299.357

# Row 139
$ws.Range("H139").Value = 72759.664
$ws.Range("J139").Value = 88139.5
$ws.Range("L139").Value = 88139.5
$ws.Range("N139").Value = -98419.5


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 124
$ws.Range("H124").Value = 40779
$ws.Range("J124").Value = 40779
$ws.Range("L124").Value = 40779
$ws.Range("N124").Value = -50599


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 87
$ws.Range("H87").Value = 28000
$ws.Range("J87").Value = 28000
$ws.Range("L87").Value = 28000
$ws.Range("N87").Value = -30246

# Row 88
$ws.Range("H88").Value = 38541.6
$ws.Range("J88").Value = 43177
$ws.Range("L88").Value = 43177
$ws.Range("N88").Value = -44033

# Row 90
$ws.Range("H90").Value = 28000
$ws.Range("J90").Value = 28000
$ws.Range("L90").Value = 84000
$ws.Range("N90").Value = -95232

# Row 91
$ws.Range("H91").Value = 38541.6
$ws.Range("J91").Value = 43177
$ws.Range("L91").Value = 43177
$ws.Range("N91").Value = -46141

# Row 93
$ws.Range("H93").Value = 1445.871
$ws.Range("I93").Value = 1110.5
$ws.Range("J93").Value = 1605.5714
$ws.Range("K93").Value = 1110.5
$ws.Range("L93").Value = 1605.5714
$ws.Range("M93").Value = 137.5
$ws.Range("N93").Value = -4101.5714

# Row 111
$ws.Range("H111").Value = 45387
$ws.Range("J111").Value = 45387
$ws.Range("L111").Value = 45387
$ws.Range("N111").Value = -53567

# Row 127
$ws.Range("H127").Value = 50705.25
$ws.Range("J127").Value = 50705.25
$ws.Range("L127").Value = 50705.25
$ws.Range("N127").Value = -60625.25

# Row 128
$ws.Range("H128").Value = 48429
$ws.Range("J128").Value = 48429
$ws.Range("L128").Value = 48429
$ws.Range("N128").Value = -58389


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 16
$ws.Range("H16").Value = 45996
$ws.Range("J16").Value = 45996
$ws.Range("L16").Value = 45996
$ws.Range("N16").Value = -46580

# Row 81
$ws.Range("H81").Value = 1228.5714
$ws.Range("I81").Value = 1460
$ws.Range("J81").Value = 650
$ws.Range("K81").Value = 2920
$ws.Range("L81").Value = 1300
$ws.Range("M81").Value = -1859
$ws.Range("N81").Value = -3422

# Row 84
$ws.Range("H84").Value = 1228.5714
$ws.Range("I84").Value = 1460
$ws.Range("J84").Value = 650
$ws.Range("K84").Value = 14600
$ws.Range("L84").Value = 6500
$ws.Range("M84").Value = -9296
$ws.Range("N84").Value = -17108

# Row 109
$ws.Range("H109").Value = 39373
$ws.Range("J109").Value = 39373
$ws.Range("L109").Value = 39373
$ws.Range("N109").Value = -42147

# Row 131
$ws.Range("H131").Value = 50578.668
$ws.Range("J131").Value = 50578.668
$ws.Range("L131").Value = 50578.668
$ws.Range("N131").Value = -60658.668

